{"js": "// Fix systematic spacing issue between header bar and body text\n// 1. Collapse the three detailed CORE COMPETENCIES paragraphs into a single\n//    summary line.\n// 2. Add a new \"TECHNICAL SKILLS\" section (heading + three summary lines)\n//    after the \"Built comprehensive survey operations platform...\" bullet.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\n// --- Step 1: collapse the CORE COMPETENCIES detail paragraphs ---------\nlet researchIdx = -1;\nlet programmingIdx = -1;\nlet dataIdx = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (researchIdx === -1 && t.indexOf(\"Research and Analytics: Survey Methodology\") === 0) {\n    researchIdx = i;\n  } else if (programmingIdx === -1 && t.indexOf(\"Programming and Development: Python\") === 0) {\n    programmingIdx = i;\n  } else if (dataIdx === -1 && t.indexOf(\"Data Infrastructure: Cloud Platforms\") === 0) {\n    dataIdx = i;\n  }\n}\n\nif (researchIdx !== -1 && programmingIdx !== -1 && dataIdx !== -1) {\n  paragraphs.items[researchIdx].insertText(\n    \"Research and Analytics \u2022 Programming and Development \u2022 Data Infrastructure\",\n    \"Replace\"\n  );\n  paragraphs.items[programmingIdx].delete();\n  paragraphs.items[dataIdx].delete();\n  await context.sync();\n}\n\n// --- Step 2: insert the new TECHNICAL SKILLS section -------------------\nconst paragraphs2 = body.paragraphs;\nparagraphs2.load(\"items,text\");\nawait context.sync();\n\nlet anchorIdx = -1;\nfor (let i = 0; i < paragraphs2.items.length; i++) {\n  if (paragraphs2.items[i].text.indexOf(\"Built comprehensive survey operations platform\") !== -1) {\n    anchorIdx = i;\n    break;\n  }\n}\n\nif (anchorIdx !== -1) {\n  const anchor = paragraphs2.items[anchorIdx];\n\n  const heading = anchor.insertParagraph(\"TECHNICAL SKILLS\", \"After\");\n  await context.sync();\n\n  const researchLine = heading.insertParagraph(\n    \"RESEARCH AND ANALYTICS Survey Methodology; Statistical Analysis; Geospatial Analysis; Data Visualization\",\n    \"After\"\n  );\n  await context.sync();\n\n  const programmingLine = researchLine.insertParagraph(\n    \"PROGRAMMING AND DEVELOPMENT Python; JVM Languages; Web Technologies; Database Languages\",\n    \"After\"\n  );\n  await context.sync();\n\n  programmingLine.insertParagraph(\n    \"DATA INFRASTRUCTURE Cloud Platforms; Big Data; Databases; Geospatial\",\n    \"After\"\n  );\n  await context.sync();\n\n  // Apply the heading style last so it does not cascade onto the\n  // paragraphs inserted after it.\n  heading.style = \"Heading 2\";\n  await context.sync();\n}\n", "ps1": "# Fix systematic spacing issue between header bar and body text\n# 1. Collapse the three detailed CORE COMPETENCIES paragraphs into a single\n#    summary line.\n# 2. Add a new \"TECHNICAL SKILLS\" section (heading + three summary lines)\n#    after the \"Built comprehensive survey operations platform...\" bullet.\n\n$d = $word.ActiveDocument\n\n# --- Step 1: collapse the CORE COMPETENCIES detail paragraphs ---------\n$researchIdx = -1\n$programmingIdx = -1\n$dataIdx = -1\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n  $t = $d.Paragraphs.Item($i).Range.Text\n  if ($researchIdx -eq -1 -and $t.StartsWith(\"Research and Analytics: Survey Methodology\")) {\n    $researchIdx = $i\n  } elseif ($programmingIdx -eq -1 -and $t.StartsWith(\"Programming and Development: Python\")) {\n    $programmingIdx = $i\n  } elseif ($dataIdx -eq -1 -and $t.StartsWith(\"Data Infrastructure: Cloud Platforms\")) {\n    $dataIdx = $i\n  }\n}\n\nif ($researchIdx -ne -1 -and $programmingIdx -ne -1 -and $dataIdx -ne -1) {\n  $researchRange = $d.Paragraphs.Item($researchIdx).Range\n  $researchRange.MoveEnd(1, -1)\n  $researchRange.Text = \"Research and Analytics \u2022 Programming and Development \u2022 Data Infrastructure\"\n\n  # Delete the Programming and Data paragraphs (Data first so the\n  # Programming paragraph's index doesn't shift before it is removed).\n  $d.Paragraphs.Item($dataIdx).Range.Delete()\n  $d.Paragraphs.Item($programmingIdx).Range.Delete()\n}\n\n# --- Step 2: insert the new TECHNICAL SKILLS section -------------------\n$anchorIdx = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n  $t = $d.Paragraphs.Item($i).Range.Text\n  if ($t.Contains(\"Built comprehensive survey operations platform\")) {\n    $anchorIdx = $i\n    break\n  }\n}\n\nif ($anchorIdx -ne -1) {\n  # Insert four empty paragraphs right after the anchor, all still in the\n  # anchor's \"Normal\" style.\n  $insertPoint = $d.Paragraphs.Item($anchorIdx).Range\n  $insertPoint.Collapse(0)\n  $insertPoint.InsertParagraphAfter()\n\n  $insertPoint2 = $d.Paragraphs.Item($anchorIdx + 1).Range\n  $insertPoint2.Collapse(0)\n  $insertPoint2.InsertParagraphAfter()\n\n  $insertPoint3 = $d.Paragraphs.Item($anchorIdx + 2).Range\n  $insertPoint3.Collapse(0)\n  $insertPoint3.InsertParagraphAfter()\n\n  $insertPoint4 = $d.Paragraphs.Item($anchorIdx + 3).Range\n  $insertPoint4.Collapse(0)\n  $insertPoint4.InsertParagraphAfter()\n\n  # Fill in the text for each new paragraph.\n  $headingRange = $d.Paragraphs.Item($anchorIdx + 1).Range\n  $headingRange.MoveEnd(1, -1)\n  $headingRange.Text = \"TECHNICAL SKILLS\"\n\n  $researchLineRange = $d.Paragraphs.Item($anchorIdx + 2).Range\n  $researchLineRange.MoveEnd(1, -1)\n  $researchLineRange.Text = \"RESEARCH AND ANALYTICS Survey Methodology; Statistical Analysis; Geospatial Analysis; Data Visualization\"\n\n  $programmingLineRange = $d.Paragraphs.Item($anchorIdx + 3).Range\n  $programmingLineRange.MoveEnd(1, -1)\n  $programmingLineRange.Text = \"PROGRAMMING AND DEVELOPMENT Python; JVM Languages; Web Technologies; Database Languages\"\n\n  $dataLineRange = $d.Paragraphs.Item($anchorIdx + 4).Range\n  $dataLineRange.MoveEnd(1, -1)\n  $dataLineRange.Text = \"DATA INFRASTRUCTURE Cloud Platforms; Big Data; Databases; Geospatial\"\n\n  # Apply the Heading 2 style last so it does not cascade onto the\n  # paragraphs inserted after it.\n  $d.Paragraphs.Item($anchorIdx + 1).Style = \"Heading 2\"\n}\n"}
